# "added data for Feb 22-23"
#
# Two new observation rows are added to the CottonObserved sheet for
# 2024-02-22 (date serial 45345):
#   - ForestHill2023IrrigationFull  -> inserted as new row 20 (pushes the
#     existing "...IrrigationPartial" rows down by one)
#   - ForestHill2023IrrigationPartial -> appended as the new last data
#     row (row 39)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row for ForestHill2023IrrigationFull (2024-02-22) -----------------
# Inserting at row 20 shifts the old row 20 (and everything below it) down
# by one, carrying formatting along exactly like a manual "Insert Row" in
# Excel.
$ws.Rows.Item(20).Insert()

$ws.Range("A20").Value = "ForestHill2023IrrigationFull"
$ws.Range("B20").Value = 45345
$ws.Range("C20").Value = 1420
$ws.Range("F20").Value = 7.135758
$ws.Range("G20").Value = 0.25569019509642193
$ws.Range("K20").NumberFormat = "0.00"
$ws.Range("K20").Value = 170.18571428571428
$ws.Range("L20").Value = 147.1142857142857

# --- New row for ForestHill2023IrrigationPartial (2024-02-22) --------------
# This is simply the next free row after the existing data (old row 37,
# now row 38 after the insert above).
$ws.Range("A39").Value = "ForestHill2023IrrigationPartial"
$ws.Range("B39").NumberFormat = "d-mmm-yy"
$ws.Range("B39").Value = 45345
$ws.Range("C39").Value = 1375
$ws.Range("F39").Value = 6.2965385000000005
$ws.Range("G39").Value = 0.3468610456650591
$ws.Range("K39").NumberFormat = "0.00"
$ws.Range("K39").Value = 126.25000000000001
$ws.Range("L39").Value = 83.516666666666666

# --- Keep the hidden _FilterDatabase name in sync with the extra row -------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "CottonObserved!_FilterDatabase") {
        $n.RefersTo = "=CottonObserved!`$A`$1:`$EQ`$2580"
    }
}

# --- Leave the selection where the user finished typing ---------------------
$ws.Range("G39").Select()
